# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# Column G (header "K") holds a recomputed per-row stat; this replaces the
# previously-stored "Strike#" derived values with the regenerated "K" values
# for each data row (rows 2-50 of the sheet; rows with an unchanged value of
# 0 are left as-is).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> new value for column G ("K")
$kValues = @{
    2  = 2
    3  = 1
    4  = 1
    5  = 2
    7  = 1
    8  = 2
    9  = 0
    10 = 1
    11 = 0
    12 = 0
    13 = 2
    14 = 1
    15 = 1
    16 = 0
    18 = 2
    19 = 1
    20 = 2
    21 = 0
    22 = 1
    23 = 2
    24 = 2
    25 = 0
    26 = 1
    27 = 0
    28 = 0
    29 = 0
    30 = 1
    31 = 1
    32 = 0
    33 = 2
    34 = 1
    35 = 1
    36 = 2
    37 = 1
    38 = 1
    39 = 0
    41 = 1
    42 = 1
    43 = 0
    44 = 0
    45 = 0
    46 = 0
    47 = 1
    48 = 1
    49 = 2
    50 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
